$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Shift header values (row 15) and comments one column to the right, from V:BS -> W:BT,
# processing right-to-left so we never overwrite a value before reading it.
# Column V gets the brand-new "culture_collection" field.

$ws.Range("BT15").Value = 'wastewater_type'
$ws.Range("BT15").AddComment('the origin of wastewater such as human waste, rainfall, storm drains, etc.')

$ws.Range("BS15").Value = 'trophic_level'
$ws.Range("BS15").Comment.Text('Feeding position in food chain (eg., chemolithotroph)')
$ws.Range("BR15").Value = 'tot_phosphate'
$ws.Range("BR15").Comment.Text('total amount or concentration of phosphate')
$ws.Range("BQ15").Value = 'tot_nitro'
$ws.Range("BQ15").Comment.Text('total nitrogen content of the sample')
$ws.Range("BP15").Value = 'tertiary_treatment'
$ws.Range("BP15").Comment.Text('the process providing a final treatment stage to raise the effluent quality before it is discharged to the receiving environment')
$ws.Range("BO15").Value = 'temperature'
$ws.Range("BO15").Comment.Text('temperature of the sample at time of sampling')
$ws.Range("BN15").Value = 'suspend_solids'
$ws.Range("BN15").Comment.Text('concentration of substances including a wide variety of material, such as silt, decaying plant and animal matter, etc,; can include multiple substances')
$ws.Range("BM15").Value = 'subspecf_gen_lin'
$ws.Range("BM15").Comment.Text('Information about the genetic distinctness of the lineage (eg., biovar, serovar)')
$ws.Range("BL15").Value = 'source_material_id'
$ws.Range("BL15").Comment.Text('unique identifier assigned to a material sample used for extracting nucleic acids, and subsequent sequencing. The identifier can refer either to the original material collected or to any derived sub-samples.')
$ws.Range("BK15").Value = 'soluble_org_mat'
$ws.Range("BK15").Comment.Text('concentration of substances such as urea, fruit sugars, soluble proteins, drugs, pharmaceuticals, etc.')
$ws.Range("BJ15").Value = 'soluble_inorg_mat'
$ws.Range("BJ15").Comment.Text('concentration of substances such as ammonia, road-salt, sea-salt, cyanide, hydrogen sulfide, thiocyanates, thiosulfates, etc.')
$ws.Range("BI15").Value = 'sodium'
$ws.Range("BI15").Comment.Text('sodium concentration')
$ws.Range("BH15").Value = 'sludge_retent_time'
$ws.Range("BH15").Comment.Text('the time activated sludge remains in reactor')
$ws.Range("BG15").Value = 'sewage_type'
$ws.Range("BG15").Comment.Text('type of wastewater treatment plant as municipial or industrial')
$ws.Range("BF15").Value = 'secondary_treatment'
$ws.Range("BF15").Comment.Text('the process for substantially degrading the biological content of the sewage')
$ws.Range("BE15").Value = 'samp_vol_we_dna_ext'
$ws.Range("BE15").Comment.Text('volume (mL) or weight (g) of sample processed for DNA extraction')
$ws.Range("BD15").Value = 'samp_store_temp'
$ws.Range("BD15").Comment.Text('temperature at which sample was stored, e.g. -80')
$ws.Range("BC15").Value = 'samp_store_loc'
$ws.Range("BC15").Comment.Text('location at which sample was stored, usually name of a specific freezer/room')
$ws.Range("BB15").Value = 'samp_store_dur'
$ws.Range("BB15").Comment.Text('duration for which sample was stored')
$ws.Range("BA15").Value = 'samp_size'
$ws.Range("BA15").Comment.Text('Amount or size of sample (volume, mass or area) that was collected')
$ws.Range("AZ15").Value = 'samp_salinity'
$ws.Range("AZ15").Comment.Text('salinity of sample, i.e. measure of total salt concentration')
$ws.Range("AY15").Value = 'samp_mat_process'
$ws.Range("AY15").Comment.Text('Processing applied to the sample during or after isolation')
$ws.Range("AX15").Value = 'samp_collect_device'
$ws.Range("AX15").Comment.Text('Method or device employed for collecting sample')
$ws.Range("AW15").Value = 'rel_to_oxygen'
$ws.Range("AW15").Comment.Text('Aerobic or anaerobic')
$ws.Range("AV15").Value = 'reactor_type'
$ws.Range("AV15").Comment.Text('anaerobic digesters can be designed and engineered to operate using a number of different process configurations, as batch or continuous, mesophilic, high solid or low solid, and single stage or multistage')
$ws.Range("AU15").Value = 'primary_treatment'
$ws.Range("AU15").Comment.Text('the process to produce both a generally homogeneous liquid capable of being treated biologically and a sludge that can be separately treated or processed')
$ws.Range("AT15").Value = 'pre_treatment'
$ws.Range("AT15").Comment.Text('the process of pre-treatment removes materials that can be easily collected from the raw wastewater')
$ws.Range("AS15").Value = 'phosphate'
$ws.Range("AS15").Comment.Text('concentration of phosphate')
$ws.Range("AR15").Value = 'ph'
$ws.Range("AR15").Comment.Text('pH measurement')
$ws.Range("AQ15").Value = 'perturbation'
$ws.Range("AQ15").Comment.Text('type of perturbation, e.g. chemical administration, physical disturbance, etc., coupled with time that perturbation occurred; can include multiple perturbation types')
$ws.Range("AP15").Value = 'pathogenicity'
$ws.Range("AP15").Comment.Text('To what is the entity pathogenic')
$ws.Range("AO15").Value = 'oxy_stat_samp'
$ws.Range("AO15").Comment.Text('oxygenation status of sample')
$ws.Range("AN15").Value = 'organism_count'
$ws.Range("AN15").Comment.Text('total count of any organism per gram or volume of sample, should include name of organism followed by count; can include multiple organism counts')
$ws.Range("AM15").Value = 'org_particles'
$ws.Range("AM15").Comment.Text('concentration of particles such as faeces, hairs, food, vomit, paper fibers, plant material, humus, etc.')
$ws.Range("AL15").Value = 'nitrate'
$ws.Range("AL15").Comment.Text('concentration of nitrate')
$ws.Range("AK15").Value = 'misc_param'
$ws.Range("AK15").Comment.Text('any other measurement performed or parameter collected, that is not listed here')
$ws.Range("AJ15").Value = 'locus_tag_prefix'
$ws.Range("AJ15").Comment.Text('A locus tag prefix required for an annotated genome, http://www.ddbj.nig.ac.jp/sub/locus_tag-e.html')
$ws.Range("AI15").Value = 'isolation_source'
$ws.Range("AI15").Comment.Text('Describes the physical, environmental and/or local geographical source of the biological sample from which the sample was derived.')
$ws.Range("AH15").Value = 'inorg_particles'
$ws.Range("AH15").Comment.Text('concentration of particles such as sand, grit, metal particles, ceramics, etc.; can include multiple particles')
$ws.Range("AG15").Value = 'indust_eff_percent'
$ws.Range("AG15").Comment.Text('percentage of industrial effluents received by wastewater treatment plant')
$ws.Range("AF15").Value = 'host_taxid'
$ws.Range("AF15").Comment.Text('NCBI taxonomy ID of the host, e.g. 9606')
$ws.Range("AE15").Value = 'host'
$ws.Range("AE15").Comment.Text('The natural (as opposed to laboratory) host to the organism from which the sample was obtained. Use the full taxonomic name, eg, "Homo sapiens".')
$ws.Range("AD15").Value = 'health_state'
$ws.Range("AD15").Comment.Text('Health or disease status of sample at time of collection')
$ws.Range("AC15").Value = 'gaseous_substances'
$ws.Range("AC15").Comment.Text('amount or concentration of substances such as hydrogen sulfide, carbon dioxide, methane, etc.; can include multiple substances')
$ws.Range("AB15").Value = 'extrachrom_elements'
$ws.Range("AB15").Comment.Text('Plasmids that have significance phenotypic consequence')
$ws.Range("AA15").Value = 'estimated_size'
$ws.Range("AA15").Comment.Text('Estimated size of genome')
$ws.Range("Z15").Value = 'encoded_traits'
$ws.Range("Z15").Comment.Text('Traits like antibiotic resistance/xenobiotic degration phenotypes/converting phage genes')
$ws.Range("Y15").Value = 'emulsions'
$ws.Range("Y15").Comment.Text('amount or concentration of substances such as paints, adhesives, mayonnaise, hair colorants, emulsified oils, etc.; can include multiple emulsion types')
$ws.Range("X15").Value = 'efficiency_percent'
$ws.Range("X15").Comment.Text('percentage of volatile solids removed from the anaerobic digestor')
$ws.Range("W15").Value = 'depth'
$ws.Range("W15").Comment.Text('Depth is defined as the vertical distance below surface, e.g. for sediment or soil samples depth is measured from sediment or soil surface, respectively. Depth can be reported as an interval for subsurface samples.')
$ws.Range("V15").Value = 'culture_collection'
$ws.Range("V15").Comment.Text('Name of source institute and unique culture identifier. See the description for the proper format and list of allowed institutes, http://www.insdc.org/controlled-vocabulary-culturecollection-qualifier')
